# New metal interference data: a "test" MetalNumber entry is added for both
# the "hollow" and "solid" MetalShape groups. The "hollow" group gains a new
# row (test) right after its existing rows, and the two previously-duplicated
# trailing "solid" rows are replaced so that "solid" also ends with a single
# new "test" row. Net effect: rows 9-17 of Sheet1 take on new contents while
# the overall used range (A1:D17) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(9,  "hollow", "test",       "1.5011101014405675E-6", "1.197085756861482E-5"),
    @(10, "solid",  "Control",    "5.7947074917335716E-7", "3.3658878039038144E-6"),
    @(11, "solid",  "LC Steel",   "4.6051638832493007E-7", "5.3854028433374753E-6"),
    @(12, "solid",  "416 SS",     "1.4205038478010832E-6", "7.6430713903476887E-6"),
    @(13, "solid",  "304 SS",     "3.2421171726750339E-6", "1.5780326137289526E-5"),
    @(14, "solid",  "6061 Al",    "5.2880230770178713E-6", "9.6326158829035823E-6"),
    @(15, "solid",  "Ti Grade 5", "2.8061680933257652E-6", "3.6577892924771489E-5"),
    @(16, "solid",  "Copper",     "1.8797195762635652E-6", "1.1754396426082914E-5"),
    @(17, "solid",  "test",       "1.5011101014405675E-6", "1.197085756861482E-5")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = [double]$row[3]
    $ws.Cells.Item($r, 4).Value = [double]$row[4]
}
